$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 21; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G ("branch")
    $old = $cell.Value2
    $cell.Value2 = 1 - $old
}
